# Add a new "TextBox 3" shape to slide 4 ("Live Tech Demo"), containing the
# five demo talking points, matching the target OOXML diff:
#
#   <p:sp> id="4" name="TextBox 3"  (appended after the existing Group 8,
#   i.e. as the 4th top-level shape on the slide)
#     spPr: off=(1444487,3154017) ext=(7733815,1477328), noFill
#     bodyPr: wrap="square" + spAutoFit
#     5 paragraphs, each a single run, colour 2C3C43, lang en-GB

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# The engine assigns each new shape's Id/Name from a per-slide monotonic
# counter (independent of the max Id already used by existing shapes).
# On this slide that counter starts at 3, but the target shape needs Id=4 /
# Name="TextBox 3". Adding (and immediately discarding) a throwaway textbox
# first burns counter value 3 so the real textbox we keep lands on 4,
# matching "TextBox 3" exactly like in the authored file.
$placeholder = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$placeholder.Delete()

# AddTextbox(Orientation, Left, Top, Width, Height) takes Left/Top/Width/
# Height in points, while the diff's xfrm is in EMU (1 pt = 12700 EMU) -
# convert so the emitted <a:off>/<a:ext> land on the exact EMU values.
$left   = 1444487 / 12700
$top    = 3154017 / 12700
$width  = 7733815 / 12700
$height = 1477328 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 3"

# <a:noFill/>
$tb.Fill.Visible = 0

$tf = $tb.TextFrame

# Set the language before the text so every run created picks it up
# (lang="en-GB" on every <a:rPr>, not just the first paragraph).
$tf.TextRange.LanguageID = "en-GB"

# One paragraph per bullet line - `r (carriage return) starts a new <a:p>.
$tf.TextRange.Text = "Show login`rShow some questions be answered`rShow the garden`rExplain the impact`rShow possibly the NLP"

# Font colour 2C3C43 on every run -> <a:solidFill><a:srgbClr val="2C3C43"/></a:solidFill>
$tf.TextRange.Font.Color.RGB = 0x433C2C

# wrap="square" + <a:spAutoFit/> on <a:bodyPr> (set after the text so the
# autofit recompute uses the final text, keeping the authored cy).
$tf.WordWrap = -1
$tf.AutoSize = 1
